$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(3, 8).Value = 4.2
$ws.Cells.Item(3, 9).Value = 3.9
$ws.Cells.Item(4, 10).Value = 1.04
$ws.Cells.Item(4, 11).Value = 13
$ws.Cells.Item(6, 14).Value = 2.02
$ws.Cells.Item(6, 15).Value = 1.88
$ws.Cells.Item(7, 10).Value = 1.08
$ws.Cells.Item(7, 11).Value = 8
$ws.Cells.Item(7, 30).Value = 900
$ws.Cells.Item(8, 14).Value = 1.57
$ws.Cells.Item(8, 15).Value = 2.35
$ws.Cells.Item(12, 8).Value = 2.95
$ws.Cells.Item(12, 9).Value = 2.4
$ws.Cells.Item(12, 14).Value = 2.37
$ws.Cells.Item(12, 15).Value = 1.45
$ws.Cells.Item(12, 18).Value = 2.02
$ws.Cells.Item(12, 19).Value = 1.62
$ws.Cells.Item(12, 20).Value = 7.1
$ws.Cells.Item(12, 21).Value = 13.5
$ws.Cells.Item(12, 23).Value = 37
$ws.Cells.Item(12, 26).Value = 6.6
$ws.Cells.Item(12, 29).Value = 120
$ws.Cells.Item(12, 34).Value = 25
$ws.Cells.Item(21, 7).Value = 1.45
$ws.Cells.Item(21, 8).Value = 4
$ws.Cells.Item(21, 9).Value = 8
$ws.Cells.Item(21, 11).Value = 8.5
$ws.Cells.Item(21, 27).Value = 8
$ws.Cells.Item(21, 28).Value = 23
$ws.Cells.Item(21, 31).Value = 13
$ws.Cells.Item(21, 36).Value = 51
$ws.Cells.Item(23, 14).Value = 2.2
$ws.Cells.Item(23, 15).Value = 1.65
$ws.Cells.Item(25, 7).Value = 3.6
$ws.Cells.Item(25, 9).Value = 2.15
$ws.Cells.Item(25, 11).Value = 9.5
$ws.Cells.Item(25, 14).Value = 1.98
$ws.Cells.Item(25, 15).Value = 1.83
$ws.Cells.Item(25, 16).Value = 1.4
$ws.Cells.Item(25, 17).Value = 2.75
$ws.Cells.Item(25, 20).Value = 11
$ws.Cells.Item(25, 26).Value = 9.5
$ws.Cells.Item(26, 7).Value = 1.6
$ws.Cells.Item(26, 8).Value = 4
$ws.Cells.Item(26, 9).Value = 5.5
$ws.Cells.Item(26, 14).Value = 1.73
$ws.Cells.Item(26, 15).Value = 2.08
$ws.Cells.Item(26, 26).Value = 13
$ws.Cells.Item(26, 30).Value = 201
$ws.Cells.Item(27, 7).Value = 3.1
$ws.Cells.Item(27, 9).Value = 2.15
$ws.Cells.Item(27, 22).Value = 11
$ws.Cells.Item(27, 24).Value = 23
$ws.Cells.Item(27, 34).Value = 21
$ws.Cells.Item(27, 36).Value = 26
$ws.Cells.Item(29, 7).Value = 1.38
$ws.Cells.Item(29, 9).Value = 7.5
$ws.Cells.Item(29, 11).Value = 17
$ws.Cells.Item(29, 28).Value = 19
$ws.Cells.Item(31, 7).Value = 2.67
$ws.Cells.Item(31, 8).Value = 3.05
$ws.Cells.Item(31, 9).Value = 2.6
$ws.Cells.Item(31, 12).Value = 1.4
$ws.Cells.Item(31, 13).Value = 2.52
$ws.Cells.Item(31, 14).Value = 2.15
$ws.Cells.Item(31, 15).Value = 1.55
$ws.Cells.Item(31, 16).Value = 1.5
$ws.Cells.Item(31, 17).Value = 2.27
$ws.Cells.Item(31, 18).Value = 1.87
$ws.Cells.Item(31, 20).Value = 7.2
$ws.Cells.Item(31, 21).Value = 12.5
$ws.Cells.Item(31, 22).Value = 10.25
$ws.Cells.Item(31, 23).Value = 30
$ws.Cells.Item(31, 24).Value = 25
$ws.Cells.Item(31, 25).Value = 40
$ws.Cells.Item(31, 26).Value = 7.5
$ws.Cells.Item(31, 27).Value = 5.9
$ws.Cells.Item(31, 31).Value = 7.1
$ws.Cells.Item(31, 32).Value = 12
$ws.Cells.Item(31, 33).Value = 10
$ws.Cells.Item(31, 34).Value = 29
$ws.Cells.Item(31, 35).Value = 24
$ws.Cells.Item(31, 36).Value = 40
$ws.Cells.Item(32, 7).Value = 4
$ws.Cells.Item(32, 8).Value = 3.05
$ws.Cells.Item(32, 9).Value = 1.95
$ws.Cells.Item(32, 10).Value = 1.12
$ws.Cells.Item(32, 11).Value = 4.55
$ws.Cells.Item(32, 12).Value = 1.57
$ws.Cells.Item(32, 13).Value = 2.12
$ws.Cells.Item(32, 14).Value = 2.6
$ws.Cells.Item(32, 15).Value = 1.38
$ws.Cells.Item(32, 16).Value = 1.6
$ws.Cells.Item(32, 17).Value = 2.07
$ws.Cells.Item(32, 18).Value = 2.32
$ws.Cells.Item(32, 19).Value = 1.47
$ws.Cells.Item(32, 20).Value = 7.8
$ws.Cells.Item(32, 21).Value = 19
$ws.Cells.Item(32, 22).Value = 15.5
$ws.Cells.Item(32, 23).Value = 65
$ws.Cells.Item(32, 24).Value = 55
$ws.Cells.Item(32, 25).Value = 90
$ws.Cells.Item(32, 26).Value = 4.9
$ws.Cells.Item(32, 27).Value = 6.4
$ws.Cells.Item(32, 28).Value = 24
$ws.Cells.Item(32, 29).Value = 200
$ws.Cells.Item(32, 31).Value = 4.9
$ws.Cells.Item(32, 32).Value = 7.4
$ws.Cells.Item(32, 33).Value = 9.75
$ws.Cells.Item(32, 34).Value = 16.5
$ws.Cells.Item(32, 35).Value = 22
$ws.Cells.Item(32, 36).Value = 55
$ws.Cells.Item(33, 7).Value = 1.8
$ws.Cells.Item(33, 8).Value = 3.45
$ws.Cells.Item(33, 9).Value = 4.5
$ws.Cells.Item(33, 18).Value = 2.1
$ws.Cells.Item(33, 19).Value = 1.65
$ws.Cells.Item(33, 20).Value = 5.6
$ws.Cells.Item(33, 21).Value = 7.8
$ws.Cells.Item(33, 23).Value = 15
$ws.Cells.Item(33, 24).Value = 17.5
$ws.Cells.Item(33, 29).Value = 150
$ws.Cells.Item(33, 31).Value = 9.75
$ws.Cells.Item(33, 32).Value = 25
$ws.Cells.Item(33, 33).Value = 16.5
$ws.Cells.Item(33, 34).Value = 90
$ws.Cells.Item(33, 35).Value = 60
$ws.Cells.Item(33, 36).Value = 75
$ws.Cells.Item(34, 11).Value = 9
$ws.Cells.Item(34, 12).Value = 1.36
$ws.Cells.Item(34, 13).Value = 3
$ws.Cells.Item(34, 28).Value = 15
$ws.Cells.Item(34, 30).Value = 301
$ws.Cells.Item(34, 32).Value = 9
$ws.Cells.Item(41, 7).Value = 3.05
$ws.Cells.Item(41, 9).Value = 2.12
$ws.Cells.Item(41, 13).Value = 4.3
$ws.Cells.Item(41, 14).Value = 1.55
$ws.Cells.Item(41, 15).Value = 2.3
$ws.Cells.Item(41, 18).Value = 1.47
$ws.Cells.Item(41, 19).Value = 2.5
$ws.Cells.Item(41, 22).Value = 10.75
$ws.Cells.Item(41, 23).Value = 40
$ws.Cells.Item(41, 24).Value = 23
$ws.Cells.Item(41, 29).Value = 35
$ws.Cells.Item(41, 31).Value = 11.5
$ws.Cells.Item(41, 32).Value = 13.5
$ws.Cells.Item(41, 34).Value = 22
$ws.Cells.Item(41, 35).Value = 14.5
$ws.Cells.Item(41, 36).Value = 18.5
$ws.Cells.Item(42, 7).Value = 3.4
$ws.Cells.Item(42, 8).Value = 3.7
$ws.Cells.Item(42, 9).Value = 1.93
$ws.Cells.Item(42, 11).Value = 8.75
$ws.Cells.Item(42, 13).Value = 4.1
$ws.Cells.Item(42, 15).Value = 2.2
$ws.Cells.Item(42, 16).Value = 1.31
$ws.Cells.Item(42, 17).Value = 3.15
$ws.Cells.Item(42, 18).Value = 1.55
$ws.Cells.Item(42, 19).Value = 2.3
$ws.Cells.Item(42, 20).Value = 13.5
$ws.Cells.Item(42, 21).Value = 21
$ws.Cells.Item(42, 22).Value = 11.75
$ws.Cells.Item(42, 24).Value = 26
$ws.Cells.Item(42, 25).Value = 27
$ws.Cells.Item(42, 26).Value = 8.75
$ws.Cells.Item(42, 27).Value = 7.4
$ws.Cells.Item(42, 28).Value = 12
$ws.Cells.Item(42, 29).Value = 45
$ws.Cells.Item(42, 31).Value = 9.5
$ws.Cells.Item(42, 32).Value = 10.75
$ws.Cells.Item(42, 34).Value = 18
$ws.Cells.Item(42, 35).Value = 14
